$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row above row 4 for the new "BardBotPath" setting (mirrors the
# existing VLCPath row) and push the existing "logF_BusinessProcessName" row
# down to row 5.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "BardBotPath"
$ws.Range("C4").Value = "Path to BardBot python bot."
$ws.Range("B4").Value = "D:\Revature\220425-UiPath\music-majors\BardBotPerformer\PerformerPythonBot\dist\BardBotP.exe"

# New row keeps the default (non-wrapped) row height like the other simple rows.
$ws.Rows.Item(4).RowHeight = 14.25

# Column B needs to widen to fit the long new path value.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 93.142857

# Update the active selection to where the user ended up editing.
$ws.Range("B12").Select()
